# Add a "fees" column to the currency_movements sheet, inserted between the
# existing "amount" (C) and "currency" (D) columns, defaulting every existing
# row's fee to 0.
# treat wire transfers only optionally as exchanges, address some further feedback

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_movements")

$lastRow = $ws.UsedRange.Rows.Count

# Shift the old columns D:E (currency, comment) one slot to the right,
# into E:F, carrying their values/styles along, and open up a blank column D.
$ws.Columns.Item(4).Insert()

# New header cell for the inserted column.
$ws.Cells.Item(1, 4).Value = "fees"

# Every existing data row gets a default fee of 0. Unlike the cells that
# Insert() shifted rightwards (which keep their original number formatting),
# the brand-new column D cells should have plain/default formatting.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = 0
    $cell.ClearFormats()
}
